# "merice plages update id, attention manque d'info"
# The "dico" sheet's Villes row (row 2) listed E2 as "identifiant" (an
# identifier/"id") for the vl_nbr_touristes field; the author corrects this
# to "obligatoire" (mandatory) instead, and leaves the selection on the
# corrected cell (E2) when saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dico")

# Correct the "Contraintes" column for the vl_nbr_touristes row.
$ws.Range("E2").Value = "obligatoire"

# Leave the selection on the cell that was just edited.
$ws.Activate()
[void]$ws.Range("E2").Select()
